# Apply the "find ems table row" changes to the EMS/small-packet tariff sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text / date updates (rich-text cells, shared strings) ---
$ws.Range("D3").Value = '04 марта 2024 № 222'
$ws.Range("H3").Value = 'от 04/03.2024 №222'
$ws.Range("D5").Value = 'в действии c 01.04.2024'
$ws.Range("H5").Value = 'в действии c 01.04.2024'

# --- Footnote wording change: "Осторожно" -> "Хрупкое" ---
$ws.Range("G20").Value = 'За хрупкие мелкие пакеты с объявленной ценностью  с отметкой "Хрупкое" взимается надбавка к плате за  пересылку  мелкого  пакета   в  размере  50%.  На  плату  за  объявленную  ценность  мелкого  пакета надбавка не распространяется'

# --- Tariff table numeric updates ---
$ws.Range("D8").Value = 2.7
$ws.Range("H8").Value = 2.25

$ws.Range("D10").Value = 3
$ws.Range("H10").Value = 2.5

$ws.Range("D11").Value = 3.48
$ws.Range("H11").Value = 3.2

$ws.Range("D15").Value = 0.84
$ws.Range("H15").Value = 0.7

$ws.Range("D16").Value = 3
$ws.Range("H16").Value = 2.5

$ws.Range("D17").Value = 3.48
$ws.Range("H17").Value = 3.2

$ws.Range("D23").Value = 1.86
$ws.Range("H23").Value = 1.55

$ws.Range("D29").Value = 3.3

$ws.Range("D32").Value = 0.8

$ws.Range("H43").Value = 2.35
$ws.Range("H44").Value = 3.1
$ws.Range("H45").Value = 3.9

$ws.Range("D47").Value = 0.8

$ws.Range("D55").Value = 2.82
$ws.Range("D56").Value = 3.72
$ws.Range("D57").Value = 9.72

$ws.Range("D80").Value = 0.78
$ws.Range("D81").Value = 2.52

$ws.Range("H82").Value = 0.65
$ws.Range("H83").Value = 2.1
